# Automatic update of files.
#
# The underlying data rows (2-22, excluding 12 and 13 which are unaffected)
# got re-paired/re-ordered against their Id (column A) - i.e. the full
# content of each row A:AY was permuted to a different row position.
# This script snapshots the original content of every affected row and then
# writes each snapshot back into its destination row according to the
# mapping derived from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns Y (Startdatum) and AA (Slutdatum) hold plain ISO date strings
# (e.g. "2023-08-31") stored as text. Excel's automatic type detection will
# otherwise silently convert such strings into real date serials when they
# are assigned programmatically, so force those columns to Text format
# first to preserve the original text representation.
$ws.Range("Y1:Y22").NumberFormat = "@"
$ws.Range("AA1:AA22").NumberFormat = "@"

# Rows that participate in the permutation (rows 12 and 13 are untouched).
$affectedRows = @(2,3,4,5,6,7,8,9,10,11,14,15,16,17,18,19,20,21,22)

# Destination row -> source row (source row's ORIGINAL content moves into
# the destination row). Derived from matching the Id (column A) values
# before/after the edit.
$rowMap = @{
    2  = 7
    3  = 21
    4  = 19
    5  = 4
    6  = 16
    7  = 9
    8  = 5
    9  = 10
    10 = 17
    11 = 8
    14 = 11
    15 = 2
    16 = 15
    17 = 18
    18 = 20
    19 = 22
    20 = 6
    21 = 3
    22 = 14
}

# Snapshot the full original row contents (A:AY) before making any changes,
# so the simultaneous/cyclic permutation does not clobber data that is
# still needed as a source for another row.
$snapshot = @{}
foreach ($r in $affectedRows) {
    $srcRange = $ws.Range("A" + $r + ":AY" + $r)
    $snapshot[$r] = $srcRange.Value()
}

# Write each row's new content from the snapshot of its mapped source row.
foreach ($destRow in $affectedRows) {
    $srcRow = $rowMap[$destRow]
    $dstRange = $ws.Range("A" + $destRow + ":AY" + $destRow)
    $dstRange.Value = $snapshot[$srcRow]
}
